$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-13 down to 7-14
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the weekly price data
$ws.Cells.Item(6,1).Value = 7
$ws.Cells.Item(6,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(6,3).Value = "Ñuble"
$ws.Cells.Item(6,4).Value = 44484
$ws.Cells.Item(6,5).Value = 16
$ws.Cells.Item(6,6).Value = 100112026
$ws.Cells.Item(6,7).Value = "Haba"
$ws.Cells.Item(6,8).Value = "Sin especificar"
$ws.Cells.Item(6,9).Value = "Primera"
$ws.Cells.Item(6,10).Value = 30
$ws.Cells.Item(6,11).Value = 8500
$ws.Cells.Item(6,12).Value = 9000
$ws.Cells.Item(6,13).Value = 8750
$ws.Cells.Item(6,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(6,15).Value = "Región del Maule"
$ws.Cells.Item(6,16).Value = 350
$ws.Cells.Item(6,17).Value = 25
$ws.Cells.Item(6,18).Value = "Hortaliza"
